# Update the "Counts" column (column B) of the "Courses Details" sheet
# with refreshed figures, per the source data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "(5,484)"
$ws.Range("B3").Value  = "(3,558)"
$ws.Range("B4").Value  = "(383)"
$ws.Range("B5").Value  = "(992)"

$ws.Range("B8").Value  = "(9,233)"
$ws.Range("B10").Value = "(4,130)"
$ws.Range("B12").Value = "(3,828)"
$ws.Range("B13").Value = "(3,783)"
$ws.Range("B14").Value = "(2,990)"
$ws.Range("B15").Value = "(2,823)"
$ws.Range("B16").Value = "(2,652)"
$ws.Range("B17").Value = "(2,533)"
$ws.Range("B18").Value = "(2,516)"
$ws.Range("B19").Value = "(2,484)"
$ws.Range("B20").Value = "(2,417)"
$ws.Range("B21").Value = "(2,401)"
$ws.Range("B22").Value = "(2,349)"
$ws.Range("B23").Value = "(2,328)"
$ws.Range("B24").Value = "(2,325)"
$ws.Range("B25").Value = "(2,307)"
$ws.Range("B26").Value = "(2,288)"
$ws.Range("B27").Value = "(2,283)"
$ws.Range("B28").Value = "(2,261)"
$ws.Range("B29").Value = "(1,927)"
$ws.Range("B30").Value = "(1,813)"
$ws.Range("B31").Value = "(1,643)"
$ws.Range("B32").Value = "(1,588)"
$ws.Range("B33").Value = "(1,353)"
$ws.Range("B34").Value = "(1,338)"
$ws.Range("B35").Value = "(1,332)"
$ws.Range("B36").Value = "(261)"

$ws.Range("B38").Value = "(38)"
